# Mise à jour de l'application
# Adds a new training-day column (DS) right after the last existing one (DR),
# copying the date/format from DR and filling in the attendance letters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy the formatting (number format / style) of column DR into the new
#    column DS so the new date column and its "P"/"B" cells look identical
#    to the existing ones.
$ws.Range("DR1:DR32").Copy()
$ws.Range("DS1:DS32").PasteSpecial(-4122)   # xlPasteFormats

# 2) New training date in the header row.
$ws.Range("DS1").Value = 46071

# 3) Attendance values for the new day, one per player row.
$values = @{
    2  = "P"
    3  = "P"
    4  = "P"
    5  = "P"
    6  = "B"
    7  = "P"
    8  = "B"
    9  = "P"
    10 = "P"
    11 = "P"
    13 = "B"
    14 = "P"
    15 = "P"
    18 = "P"
    19 = "P"
    20 = "P"
    22 = "P"
    24 = "P"
    26 = "P"
    27 = "P"
    28 = "P"
    29 = "P"
    30 = "P"
    31 = "P"
    32 = "P"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 123).Value = $values[$row]
}

# Rows 16, 17 and 25 only have the blank formatted cell (no attendance that
# day yet) and rows 12, 21, 23 stop earlier and should not get a DS cell at
# all, so explicitly clear those three to drop them from the sheet.
$ws.Range("DS12").Clear()
$ws.Range("DS21").Clear()
$ws.Range("DS23").Clear()

# 4) Update the frozen pane / current selection to reflect the newly added
#    column, matching the author's latest view state.
$ws.Range("DT11").Select()
$ws.Application.ActiveWindow.FreezePanes = $false
$ws.Range("DP1").Select()
$excel.ActiveWindow.SplitColumn = 1
$excel.ActiveWindow.SplitRow = 0
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("DT11").Select()

# 5) Recalculate so the COUNTA/COUNTIF summary columns (B, C, F, ...) pick
#    up the newly added attendance cells.
$excel.Calculate()
